$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "https://interop.esante.gouv.fr/ig/fhir/tde/StructureDefinition/EyeColor"
$meta.Range("B3").Value = "2.0.0"
$meta.Range("B8").Value = "2026-01-15T15:25:40+00:00"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("Z6").Value = "https://interop.esante.gouv.fr/ig/fhir/tde/ValueSet/EyeColorVS"

# Column Z width shrinks because the new URL text is shorter (auto-fit effect).
# Target stored width is 49.4453125; the COM layer quantizes ColumnWidth to
# ~1/6-character steps, so 48.6 is the closest input that lands on the nearest
# achievable stored value (49.5).
$elements.Columns.Item(26).ColumnWidth = 48.6
